# Actualización automática de scrims_actualizado.xlsx (2025-07-27 01:27:13)
#
# Appends new scrim-result rows to five worksheets. Each new row is built by
# cloning the row format immediately above it (so fills/borders/fonts match
# the table's existing look), writing the real values into every column, and
# - when the "Ganador" (G) column's outcome differs from the template row -
# re-coloring just that cell so "Equipo 1" / "Equipo 2" / "Empate" keep their
# normal blue / pink / gray highlighting.

$wb = $excel.ActiveWorkbook

# Blue / pink / gray fills (OLE BGR ints) used for the "Ganador" column.
$colorEquipo1 = 16770508   # RGB CCE5FF
$colorEquipo2 = 13421812   # RGB F4CCCC
$colorEmpate  = 14277081   # RGB D9D9D9

function Add-ScrimRow {
    # Positional args (named binding is unreliable in this host):
    #   1 SheetName, 2 TemplateRow, 3 NewRow, 4 Values (14 items, cols A..N)
    param(
        [string]$SheetName,
        [int]$TemplateRow,
        [int]$NewRow,
        [string[]]$Values
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Clone formatting (fills/borders/fonts) plus values from the row above,
    # then we'll overwrite the values with the real data below.
    $srcRange = $ws.Range("A" + $TemplateRow + ":N" + $TemplateRow)
    $dstRange = $ws.Range("A" + $NewRow + ":N" + $NewRow)
    $srcRange.Copy($dstRange)

    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = $cols[$i] + $NewRow
        $ws.Range($addr).Value = $Values[$i]
    }

    # Fix the "Ganador" cell's highlight if this row's outcome differs from
    # the template row's outcome (blue = Equipo 1, pink = Equipo 2, gray = Empate).
    $templateGanador = $ws.Range("G" + $TemplateRow).Value
    $newGanador = $Values[6]
    if ($newGanador -ne $templateGanador) {
        $gCell = $ws.Range("G" + $NewRow)
        if ($newGanador -eq "Equipo 1") {
            $gCell.Interior.Color = $colorEquipo1
        } elseif ($newGanador -eq "Equipo 2") {
            $gCell.Interior.Color = $colorEquipo2
        } else {
            $gCell.Interior.Color = $colorEmpate
        }
        $gCell.Font.Bold = $true
        $gCell.Borders.LineStyle = 1
    }
}

# ---------------------------------------------------------------------------
# Triple Dribble: A3:N98 -> A3:N100 (rows 99-100)
# ---------------------------------------------------------------------------
Add-ScrimRow "Triple Dribble" 98 99 @(
    "KAZE","CHARLIE","HANK","JAE-YONG","BARLEY","FRANK","Equipo 1",
    "FX|CaueBr","FX|Wesley","FX|REI DO FUT","Satisfyer🌺","Bicho🐦‍🔥","Fv7🐦‍🔥",
    "20250726T231156.000Z"
)
Add-ScrimRow "Triple Dribble" 99 100 @(
    "KAZE","CHARLIE","HANK","JAE-YONG","BARLEY","FRANK","Equipo 1",
    "FX|CaueBr","FX|Wesley","FX|REI DO FUT","Satisfyer🌺","Bicho🐦‍🔥","Fv7🐦‍🔥",
    "20250726T230911.000Z"
)

# ---------------------------------------------------------------------------
# Sneaky Fields: A3:N22 -> A3:N25 (rows 23-25)
# ---------------------------------------------------------------------------
Add-ScrimRow "Sneaky Fields" 22 23 @(
    "BEA","BULL","RICO","CORDELIUS","GUS","DRACO","Equipo 2",
    "BC*|Derrp","BC*|Jubileubr","BC*|Loko","SKC|Rhz","SKC|Kr","SKC|Prozy",
    "20250726T232131.000Z"
)
Add-ScrimRow "Sneaky Fields" 23 24 @(
    "BEA","BULL","RICO","CORDELIUS","GUS","DRACO","Equipo 1",
    "BC*|Derrp","BC*|Jubileubr","BC*|Loko","SKC|Rhz","SKC|Kr","SKC|Prozy",
    "20250726T231840.000Z"
)
Add-ScrimRow "Sneaky Fields" 24 25 @(
    "BEA","BULL","RICO","CORDELIUS","GUS","DRACO","Equipo 2",
    "BC*|Derrp","BC*|Jubileubr","BC*|Loko","SKC|Rhz","SKC|Kr","SKC|Prozy",
    "20250726T231633.000Z"
)

# ---------------------------------------------------------------------------
# Dueling Beetles: A3:N25 -> A3:N28 (rows 26-28)
# ---------------------------------------------------------------------------
Add-ScrimRow "Dueling Beetles" 25 26 @(
    "GRAY","CHARLIE","BARLEY","LUMI","AMBER","ASH","Equipo 1",
    "CASA|Pekka","Doritos🐉","CASA|Mohtep","Dreww :D","Fire Murilø:D🧸","Buk :D",
    "20250726T231751.000Z"
)
Add-ScrimRow "Dueling Beetles" 26 27 @(
    "GRAY","CHARLIE","BARLEY","LUMI","AMBER","ASH","Equipo 1",
    "CASA|Pekka","Doritos🐉","CASA|Mohtep","Dreww :D","Fire Murilø:D🧸","Buk :D",
    "20250726T231521.000Z"
)
Add-ScrimRow "Dueling Beetles" 27 28 @(
    "GRAY","CHARLIE","BARLEY","LUMI","AMBER","ASH","Equipo 2",
    "CASA|Pekka","Doritos🐉","CASA|Mohtep","Dreww :D","Fire Murilø:D🧸","Buk :D",
    "20250726T231234.000Z"
)

# ---------------------------------------------------------------------------
# Hot Potato: A3:N102 -> A3:N105 (rows 103-105)
# ---------------------------------------------------------------------------
Add-ScrimRow "Hot Potato" 102 103 @(
    "MELODIE","CORDELIUS","HANK","NITA","BERRY","MICO","Equipo 1",
    "FX|CaueBr","FX|Wesley","FX|REI DO FUT","Satisfyer🌺","Bicho🐦‍🔥","Fv7🐦‍🔥",
    "20250726T232138.000Z"
)
Add-ScrimRow "Hot Potato" 103 104 @(
    "MELODIE","CORDELIUS","HANK","NITA","BERRY","MICO","Equipo 1",
    "FX|CaueBr","FX|Wesley","FX|REI DO FUT","Satisfyer🌺","Bicho🐦‍🔥","Fv7🐦‍🔥",
    "20250726T232026.000Z"
)
Add-ScrimRow "Hot Potato" 104 105 @(
    "MELODIE","CORDELIUS","HANK","NITA","BERRY","MICO","Equipo 2",
    "FX|CaueBr","FX|Wesley","FX|REI DO FUT","Satisfyer🌺","Bicho🐦‍🔥","Fv7🐦‍🔥",
    "20250726T231826.000Z"
)

# ---------------------------------------------------------------------------
# Ring of Fire: A3:N92 -> A3:N94 (rows 93-94)
# ---------------------------------------------------------------------------
Add-ScrimRow "Ring of Fire" 92 93 @(
    "DOUG","HANK","MEEPLE","GUS","KIT","BULL","Equipo 1",
    "BC*|Derrp","BC*|Jubileubr","BC*|Loko","SKC|Rhz","SKC|Kr","SKC|Prozy",
    "20250726T231102.000Z"
)
Add-ScrimRow "Ring of Fire" 93 94 @(
    "DOUG","HANK","MEEPLE","GUS","KIT","BULL","Equipo 1",
    "BC*|Derrp","BC*|Jubileubr","BC*|Loko","SKC|Rhz","SKC|Kr","SKC|Prozy",
    "20250726T230912.000Z"
)
